# Append new "Injured" period activity rows (91-112) to the Activity sheet,
# recording a stretch of rest/light activity following an injury.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

# --- Seed the four brand-new Activity labels in the exact order they first
# appear in the final shared-strings table, so de-duplication lines up with
# the authored workbook. (Cardio is actually entered on row 108, but in the
# saved file's dedup table it sorts after "Walking " from row 111 -
# matching the author's real edit order.)
$ws.Cells.Item(96, 2).Value = "Stretching"
$ws.Cells.Item(98, 2).Value = "Walking"
$ws.Cells.Item(111, 2).Value = "Walking "
$ws.Cells.Item(108, 2).Value = "Cardio"

# Row data: Date(serial), Activity, Details, Duration(min), Time(day fraction), Intensity
$rows = @(
    @{ R=91;  Date=45285; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=92;  Date=45286; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=93;  Date=45287; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=94;  Date=45288; Activity="BW Exercises"; Details="Injured"; Duration=10; Time=0.54166666666666663; Intensity=$null },
    @{ R=95;  Date=45289; Activity="BW Exercises"; Details="Injured"; Duration=10; Time=0.58333333333333337; Intensity=$null },
    @{ R=96;  Date=45290; Activity="Stretching";   Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=97;  Date=45291; Activity="BW Exercises"; Details="Injured"; Duration=15; Time=0.625;                Intensity=$null },
    @{ R=98;  Date=45292; Activity="Walking";      Details="Injured"; Duration=25; Time=0.75;                 Intensity=$null },
    @{ R=99;  Date=45293; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=100; Date=45294; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=101; Date=45295; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=102; Date=45296; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=103; Date=45297; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=104; Date=45298; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null },
    @{ R=105; Date=45299; Activity="Walking";      Details="Injured"; Duration=30; Time=0.75;                 Intensity=7 },
    @{ R=106; Date=45300; Activity="Walking";      Details="Injured"; Duration=30; Time=0.72916666666666663; Intensity=7 },
    @{ R=107; Date=45301; Activity="Walking";      Details="Injured"; Duration=30; Time=0.875;                Intensity=7 },
    @{ R=108; Date=45302; Activity="Cardio";       Details="Injured"; Duration=15; Time=0.58333333333333337; Intensity=8 },
    @{ R=109; Date=45303; Activity="Walking";      Details="Injured"; Duration=30; Time=0.79166666666666663; Intensity=8 },
    @{ R=110; Date=45304; Activity="Cardio";       Details="Injured"; Duration=15; Time=0.54166666666666663; Intensity=7 },
    @{ R=111; Date=45305; Activity="Walking ";     Details="Injured"; Duration=30; Time=0.79166666666666663; Intensity=9 },
    @{ R=112; Date=45306; Activity="Rest";         Details="Injured"; Duration=0;  Time=$null;               Intensity=$null }
)

# Use the existing formatted rows (A90 date style, F90 time style) as format
# donors so new cells pick up identical number formats/styles instead of
# Excel inventing new custom number-format entries.
$dateFmt = $ws.Range("A90")
$timeFmt = $ws.Range("F90")

foreach ($row in $rows) {
    $r = $row.R

    $dateFmt.Copy()
    $ws.Range("A$r").PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $row.Date

    $ws.Cells.Item($r, 2).Value = $row.Activity
    $ws.Cells.Item($r, 3).Value = $row.Details
    $ws.Cells.Item($r, 4).Value = $row.Duration

    if ($null -ne $row.Time) {
        $timeFmt.Copy()
        $ws.Range("F$r").PasteSpecial(-4122)
        $ws.Cells.Item($r, 6).Value = $row.Time
    }
    if ($null -ne $row.Intensity) {
        $ws.Cells.Item($r, 7).Value = $row.Intensity
    }
}

$excel.CutCopyMode = 0

# Update sheet view: scroll position, zoom and active selection to match
# where the author ended up after entering this data.
$win = $excel.ActiveWindow
$win.ScrollRow = 96
$win.Zoom = 85
$ws.Range("D106").Select()
